$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "System coverage" rows (C9:C16, the "Tested in C99" column for the
# exit-code related checks) as written: copy the "Test written" look (green
# fill / border) from an already-completed cell (C4) and apply the text.
$ws.Range("C4").Copy()
[void]$ws.Range("C9:C16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C9:C16").Value = "Test written"

# Restore the selection Excel shows after making this edit.
[void]$ws.Range("B35").Select()
